$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (|S*|/n)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: summary statistics with labels in column A and values in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold/size-12/vertically-centered format once on a scratch cell,
# then copy it onto the summary value cells in a single formatting op so we
# don't leave unused intermediate styles behind in styles.xml.
$tpl = $ws.Range("D100")
$tpl.Font.Bold = $true
$tpl.Font.Size = 12
$tpl.VerticalAlignment = -4108

$tpl.Copy()
$ws.Range("B14:B17").PasteSpecial(-4122)
$tpl.Clear()

# Match the authored selection state
$ws.Range("A14:B17").Select()

$wb.Save()
